# adding scapula trait terms
#
# Populates the "in oba" / "in fovt" / "Y" columns (A/B) on the "trait"
# sheet for the scapula rows that are already mapped to FOVT/OBA terms,
# matches the view/selection state left behind by the edit, and restores
# a couple of cosmetic sheet settings (page orientation, axis column A
# width) that changed alongside the new data.

$wb = $excel.ActiveWorkbook

$trait = $wb.Worksheets.Item("trait")
$axis  = $wb.Worksheets.Item("axis")
$structures = $wb.Worksheets.Item("structures")

# --- New term columns on the "trait" sheet -------------------------------
# Row 3 only gets the "in oba" marker (column A).
$trait.Cells.Item(3, 1).Value = "in oba"

# Rows 4,5,6,8,9,10 get "in fovt" (A) + "Y" (B). "Y" is written first so
# the shared-string table picks up the same ordering as the authored file
# (in oba=60, Y=61, in fovt=62).
$fovtRows = 4,5,6,8,9,10
foreach ($r in $fovtRows) {
    $trait.Cells.Item($r, 2).Value = "Y"
    $trait.Cells.Item($r, 1).Value = "in fovt"
}

# --- Page setup tweak on "trait" -----------------------------------------
$trait.PageSetup.Orientation = 1   # xlPortrait

# --- Column A best-fit width on "axis" ------------------------------------
$axis.Columns.Item(1).AutoFit()

# --- Selections / active sheet, restored in the order that leaves
#     "trait" as the final active tab (matches tabSelected="1" there and
#     no tabSelected elsewhere). The "AB" sheet's view is untouched by the
#     original edit, so it is deliberately left alone here. -------------
[void]$axis.Range("B7").Select()
[void]$structures.Range("E11").Select()

$trait.Activate()
[void]$trait.Range("A11").Select()
